$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: for D-column values that look like plain numbers, force them
# to be stored as text (matching the source data which uses text-formatted
# prices), then reset the cell style back to Normal so no stray number format
# is left applied to the cell.

# Row 27/28 swap: PancakeSwap and RenderToken swap positions
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.30%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "

# Price/Volume updates for other rows
$ws.Range("D2").Value = "66.095.30"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.180.16"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "3.178.06"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.513"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("D15").Value = "3.701.02"
$ws.Range("E15").Value = "  +1.75%  "
$ws.Range("D16").Value = "66.077.18"
$ws.Range("E16").Value = "  +2.10%  "
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "3.184.11"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.88%  "
$ws.Range("E22").Value = "  +3.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +14.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("E33").Value = "  +3.19%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "487.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0893"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.122"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.298"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "0.0₃0652"
$ws.Range("E44").Value = "  +14.96%  "
$ws.Range("D45").Value = "2.891.39"
$ws.Range("E45").Value = "  -4.21%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.117"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.87%  "
$ws.Range("E51").Value = "  +3.91%  "
